$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KATALYST progress updated from 88 to 75 (row 5, column B)
$ws.Range("B5").Value = 75

# HORIZON OLE progress value cleared (row 12, column B)
$ws.Range("B12").ClearContents()
